$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "63.228.35"
$ws.Range("E2").Value = "  -1.45%  "

$ws.Range("D3").Value = "3.088.61"
$ws.Range("E3").Value = "  +0.37%  "

$ws.Range("E4").Value = "  +0.03%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "554.33"
$ws.Range("E5").Value = "  +0.31%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "137.14"
$ws.Range("E6").Value = "  -3.86%  "

$ws.Range("E7").Value = "  -0.02%  "

$ws.Range("D8").Value = "3.077.22"
$ws.Range("E8").Value = "  +0.29%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.495"
$ws.Range("E9").Value = "  +1.25%  "

$ws.Range("B10").Value = "Dogecoin"
$ws.Range("C10").Value = "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.160"
$ws.Range("E10").Value = "  +5.78%  "

$ws.Range("B11").Value = "Toncoin"
$ws.Range("C11").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "6.61"
$ws.Range("E11").Value = "  +1.93%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.454"
$ws.Range("E12").Value = "  +1.05%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "35.04"
$ws.Range("E13").Value = "  -1.51%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.0000217"
$ws.Range("E14").Value = "  +0.70%  "

$ws.Range("D15").Value = "3.580.53"
$ws.Range("E15").Value = "  +0.31%  "

$ws.Range("D16").Value = "63.270.78"
$ws.Range("E16").Value = "  -1.44%  "

$ws.Range("E17").Value = "  +0.12%  "

$ws.Range("D18").Value = "3.093.34"
$ws.Range("E18").Value = "  +0.54%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "500.92"
$ws.Range("E19").Value = "  +2.60%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.65"
$ws.Range("E20").Value = "  +1.13%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.53"
$ws.Range("E21").Value = "  -0.20%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.705"
$ws.Range("E22").Value = "  +3.87%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.27"
$ws.Range("E23").Value = "  +1.61%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "78.03"
$ws.Range("E24").Value = "  +0.93%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "12.31"
$ws.Range("E25").Value = "  -0.55%  "

$ws.Range("E26").Value = "  +0.08%  "

$ws.Range("E27").Value = "  +1.76%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "8.15"
$ws.Range("E28").Value = "  -1.18%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.02"
$ws.Range("E29").Value = "  -1.89%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.00"
$ws.Range("E30").Value = "  +0.05%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "26.22"
$ws.Range("E31").Value = "  +2.11%  "

$ws.Range("E32").Value = "  -4.61%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.11"
$ws.Range("E33").Value = "  -1.77%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "59.05"
$ws.Range("E34").Value = "  +13.37%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "532.17"
$ws.Range("E35").Value = "  -8.41%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.87"
$ws.Range("E36").Value = "  +0.25%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "5.14"

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0412"
$ws.Range("E38").Value = "  +3.13%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0792"
$ws.Range("E39").Value = "  +0.78%  "

$ws.Range("D40").Value = "3.058.21"
$ws.Range("E40").Value = "  +2.10%  "

$ws.Range("E41").Value = "  +2.01%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "8.06"
$ws.Range("E42").Value = "  -0.84%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.64"
$ws.Range("E43").Value = "  -6.05%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.253"
$ws.Range("E44").Value = "  +4.05%  "

$ws.Range("E45").Value = "  +0.02%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.06"
$ws.Range("E46").Value = "  -1.23%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "120.26"
$ws.Range("E47").Value = "  +1.48%  "

$ws.Range("B48").Value = "Stellar"
$ws.Range("C48").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.106"
$ws.Range("E48").Value = "  -0.28%  "

$ws.Range("B49").Value = "InjectiveProtocol"
$ws.Range("C49").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "23.84"
$ws.Range("E49").Value = "  -4.65%  "

$ws.Range("D50").Value = "0.0₃0496"
$ws.Range("E50").Value = "  -5.06%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.35"
$ws.Range("E51").Value = "  +68.18%  "
